$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -10.30360000000001
$ws.Range("C18").Value = -11.80309999999999
$ws.Range("C20").Value = -12.37349999999999
$ws.Range("C27").Value = -12.38869999999999
$ws.Range("C69").Value = -11.1422
$ws.Range("C76").Value = -12.5783
$ws.Range("C82").Value = -11.9004
